$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# Status column updates (all rows referencing the old "Ready for handoff" text)
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# Latest Handback DateTime column (H)
$zhcn.Range("H2").Value = "2016-03-19 12:24:43"
$zhcn.Range("H3").Value = "2016-03-19 12:24:43"

$dede.Range("H2").Value = "2016-03-19 12:24:49"
$dede.Range("H3").Value = "2016-03-19 12:24:49"

# New "Latest Target File" (F) and "Latest Handback File" (G) columns
$zhcn.Range("F2").Value = "a.md"
$zhcn.Range("F2").Style = "Hyperlink"
$zhcn.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("G2").Style = "Hyperlink"

$zhcn.Range("F3").Value = "a.md"
$zhcn.Range("F3").Style = "Hyperlink"
$zhcn.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("G3").Style = "Hyperlink"

$dede.Range("F2").Value = "a.md"
$dede.Range("F2").Style = "Hyperlink"
$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("G2").Style = "Hyperlink"

$dede.Range("F3").Value = "a.md"
$dede.Range("F3").Style = "Hyperlink"
$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("G3").Style = "Hyperlink"
